$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.272.63'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.863.21'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4708'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2906'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06540'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.84'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07937'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.84'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").Value = '1.863.07'
$ws.Range("E13").Value = '  -0.45%  '
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6805'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '264.15'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -6.00%  '
$ws.Range("D17").Value = '30.256.28'
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.76'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +8.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007449'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.17%  '
$ws.Range("D21").Value = '2.104.09'
$ws.Range("E21").Value = '  -0.57%  '
$ws.Range("E22").Value = '  -0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.267'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.172'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.32'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.189'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.90'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.951'
$ws.Range("D28").ClearFormats()
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.395'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09851'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.352'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("E32").Value = '  -0.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.029'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.01%  '
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.129'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7000'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.707'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("E38").Value = '  +0.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.623'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.353'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.01'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.944'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8432'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9997'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4154'
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.28'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.160'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '946.97'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.209'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.17'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05663'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.57%  '
